# Applies the changes described by the diff:
#  1. Remove the "_GoBack" bookmark from the "Need to create a service..." paragraph
#     (Front End Angular Changes section).
#  2. Add a "_GoBack" bookmark at the end of the "...for adding parent nodes." paragraph
#     (right before its paragraph mark).
#  3. Insert a new paragraph ("*Switch to using OpenCpu...") right after the
#     "*We need to add a function..." paragraph, before the existing blank paragraph.
#  4. Insert two new paragraphs ("Before we create a remote repo on Mordor..." and a
#     single-space paragraph) right after the existing blank paragraph, at the end of
#     the document body (before the sectPr).

$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: re-add the "_GoBack" bookmark at the end of the target paragraph -----
# Locate the paragraph that ends with "for adding parent nodes."
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*for adding parent nodes.*") {
        $targetPara = $candidate
        break
    }
}

$r = $targetPara.Range
$r.MoveEnd(1, -1)          # exclude the paragraph mark
$r.Collapse(0)             # collapse to the end (wdCollapseEnd)

# Placing a bookmark with Bookmarks.Add directly at this exact boundary position can
# misplace it, so nudge past the boundary first: insert a temporary character, collapse
# to just before it, add the bookmark there, then remove the temporary character.
$r.InsertAfter("X")
$r.Collapse(1)             # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $r)
$tempCharRange = $d.Range($r.Start, $r.Start + 1)
$tempCharRange.Delete()

# --- Step 3/4: insert the new paragraphs -------------------------------------------

# Find the "*We need to add a function..." paragraph and the trailing blank paragraph
# that immediately follows it.
$redrawPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*redraw the graph*") {
        $redrawPara = $candidate
        break
    }
}
$redrawParaIndex = $redrawPara.Index
$blankPara = $d.Paragraphs.Item($redrawParaIndex + 1)

$openCpuXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 3) New paragraph right after the "redraw the graph" paragraph (i.e. right before the
#    blank paragraph).
$switchParaXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document $openCpuXmlNs><w:body><w:p><w:r><w:t xml:space="preserve">*Switch to using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OpenCpu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> instead of R script and see if that has any impact in performance. I am hoping that it will allow us to keep object sin memory so that we don&#8217;t have to keep reading in correlation matrices.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$insertRange = $blankPara.Range
$insertRange.Collapse(1)   # wdCollapseStart -- start of the blank paragraph
$insertRange.InsertXML($switchParaXml)

# Re-resolve the blank paragraph (it may have shifted index after the insert above).
$blankPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq [char]13) {
        $precedingIndex = $candidate.Index - 1
        if ($precedingIndex -ge 1) {
            $precedingText = $d.Paragraphs.Item($precedingIndex).Range.Text
            if ($precedingText -like "*correlation matrices.*") {
                $blankPara = $candidate
                break
            }
        }
    }
}

# 4) New paragraphs right after the (still) blank paragraph, at the very end of the
#    document body.
$mordorParasXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document $openCpuXmlNs><w:body><w:p><w:r><w:t xml:space="preserve">Before we create a remote repo on Mordor, let&#8217;s first transition to using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OpenCpu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Unfortunately, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OpenCpu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is actually a server. We don&#8217;t seem to have much control over where that server exists and which directory is its root directory. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$insertRange2 = $blankPara.Range
$insertRange2.Collapse(0)  # wdCollapseEnd -- end of the blank paragraph (before its mark)
$insertRange2.InsertXML($mordorParasXml)
